$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.627.08"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "3.446.01"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'579.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'148.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'8.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.94%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("D12").Value = "4.035.96"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "'28.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.15%  "
$ws.Range("D15").Value = "3.438.94"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "62.686.70"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "'6.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'14.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").Value = "'9.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("D21").Value = "'386.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").Value = "'75.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "3.581.77"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").Value = "'0.0000114"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "'0.183"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "'7.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").Value = "'7.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -5.04%  "
$ws.Range("D34").Value = "'23.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("D35").Value = "'1.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.99%  "
$ws.Range("D36").Value = "'5.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").Value = "'31.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "'6.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").Value = "'169.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "3.480.68"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "'0.0772"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.785"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").Value = "'42.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").Value = "2.574.98"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "'6.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "'2.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("D50").Value = "'22.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("E51").Value = "  +0.02%  "
